$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-30 (A,B,C) with new data values; styles already present so plain Value assignment keeps them
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 44347
$ws.Range("C2").Value = 1.173958831738437

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 44348
$ws.Range("C3").Value = 1.078701790633609

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 44349
$ws.Range("C4").Value = 1.00952582814624

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 44350
$ws.Range("C5").Value = 1.0660650727329

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 44351
$ws.Range("C6").Value = 1.125101092896175

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 44352
$ws.Range("C7").Value = 1.002557203002014

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 44353
$ws.Range("C8").Value = 1.055794451450189

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 44354
$ws.Range("C9").Value = 1.088159681033346

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 44355
$ws.Range("C10").Value = 1.076602097672894

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 44356
$ws.Range("C11").Value = 1.073770464304884

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 44357
$ws.Range("C12").Value = 1.13271992880126

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 44358
$ws.Range("C13").Value = 1.135222861250899

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 44359
$ws.Range("C14").Value = 0.9041353125612385

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 44360
$ws.Range("C15").Value = 1.039103831891224

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 44361
$ws.Range("C16").Value = 1.24174882629108

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 44362
$ws.Range("C17").Value = 1.207198000981114

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 44363
$ws.Range("C18").Value = 1.221758621788742

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 44364
$ws.Range("C19").Value = 1.233141511266511

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 44365
$ws.Range("C20").Value = 1.202246946879449

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 44366
$ws.Range("C21").Value = 1.119892361111111

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 44367
$ws.Range("C22").Value = 1.071186713191024

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 44368
$ws.Range("C23").Value = 1.289848046068837

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 44369
$ws.Range("C24").Value = 1.249586786604676

$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 44370
$ws.Range("C25").Value = 1.167303615063964

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 44371
$ws.Range("C26").Value = 1.157776364664083

$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 44372
$ws.Range("C27").Value = 1.121360982311523

$ws.Range("A28").Value = 26
$ws.Range("B28").Value = 44373
$ws.Range("C28").Value = 1.035854541692531

$ws.Range("A29").Value = 27
$ws.Range("B29").Value = 44374
$ws.Range("C29").Value = 1.102360217934977

$ws.Range("A30").Value = 28
$ws.Range("B30").Value = 44375
$ws.Range("C30").Value = 1.280216749350342

# Extend formatting down to new rows by copying row 30 formatting (style s="1"/s="2") to rows 31-62
$ws.Range("A30:C30").Copy() | Out-Null
$ws.Range("A31:C62").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill new rows 31-62 with data values
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 44376
$ws.Range("C31").Value = 1.188741780475468

$ws.Range("A32").Value = 30
$ws.Range("B32").Value = 44377
$ws.Range("C32").Value = 1.034486146451795

$ws.Range("A33").Value = 31
$ws.Range("B33").Value = 44378
$ws.Range("C33").Value = 0.9908808254749214

$ws.Range("A34").Value = 32
$ws.Range("B34").Value = 44379
$ws.Range("C34").Value = 1.03166557591623

$ws.Range("A35").Value = 33
$ws.Range("B35").Value = 44380
$ws.Range("C35").Value = 1.087660328113348

$ws.Range("A36").Value = 34
$ws.Range("B36").Value = 44381
$ws.Range("C36").Value = 0.8968446026097271

$ws.Range("A37").Value = 35
$ws.Range("B37").Value = 44382
$ws.Range("C37").Value = 1.030911111111111

$ws.Range("A38").Value = 36
$ws.Range("B38").Value = 44383
$ws.Range("C38").Value = 0.9868629745189808

$ws.Range("A39").Value = 37
$ws.Range("B39").Value = 44384
$ws.Range("C39").Value = 0.8768915128483532

$ws.Range("A40").Value = 38
$ws.Range("B40").Value = 44385
$ws.Range("C40").Value = 0.9378790279978869

$ws.Range("A41").Value = 39
$ws.Range("B41").Value = 44386
$ws.Range("C41").Value = 0.9539435048155979

$ws.Range("A42").Value = 40
$ws.Range("B42").Value = 44387
$ws.Range("C42").Value = 0.88443538647343

$ws.Range("A43").Value = 41
$ws.Range("B43").Value = 44388
$ws.Range("C43").Value = 1.067987847222222

$ws.Range("A44").Value = 42
$ws.Range("B44").Value = 44389
$ws.Range("C44").Value = 0.9239915458937199

$ws.Range("A45").Value = 43
$ws.Range("B45").Value = 44390
$ws.Range("C45").Value = 0.9855636070853462

$ws.Range("A46").Value = 44
$ws.Range("B46").Value = 44391
$ws.Range("C46").Value = 0.8948792557403009

$ws.Range("A47").Value = 45
$ws.Range("B47").Value = 44392
$ws.Range("C47").Value = 0.8796068131168417

$ws.Range("A48").Value = 46
$ws.Range("B48").Value = 44393
$ws.Range("C48").Value = 0.8003962264150943

$ws.Range("A49").Value = 47
$ws.Range("B49").Value = 44394
$ws.Range("C49").Value = 0.7912763409961686

$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 44395
$ws.Range("C50").Value = 0.8446935626102293

$ws.Range("A51").Value = 49
$ws.Range("B51").Value = 44396
$ws.Range("C51").Value = 0.7717011558538404

$ws.Range("A52").Value = 50
$ws.Range("B52").Value = 44397
$ws.Range("C52").Value = 0.9096218572587186

$ws.Range("A53").Value = 51
$ws.Range("B53").Value = 44398
$ws.Range("C53").Value = 0.9291658929124109

$ws.Range("A54").Value = 52
$ws.Range("B54").Value = 44399
$ws.Range("C54").Value = 0.8888506760728982

$ws.Range("A55").Value = 53
$ws.Range("B55").Value = 44400
$ws.Range("C55").Value = 0.8062001318826245

$ws.Range("A56").Value = 54
$ws.Range("B56").Value = 44401
$ws.Range("C56").Value = 0.52578125

$ws.Range("A57").Value = 55
$ws.Range("B57").Value = 44402
$ws.Range("C57").Value = 0.7691309987029831

$ws.Range("A58").Value = 56
$ws.Range("B58").Value = 44403
$ws.Range("C58").Value = 0.8770076377523186

$ws.Range("A59").Value = 57
$ws.Range("B59").Value = 44404
$ws.Range("C59").Value = 0.9851051365611422

$ws.Range("A60").Value = 58
$ws.Range("B60").Value = 44405
$ws.Range("C60").Value = 0.9348814864662779

$ws.Range("A61").Value = 59
$ws.Range("B61").Value = 44406
$ws.Range("C61").Value = 1.012321344616831

$ws.Range("A62").Value = 60
$ws.Range("B62").Value = 44407
$ws.Range("C62").Value = 0.9185994363929147

